# Update the build/version string across the workbook.
#
# Old version string: "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# New version string: "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: $newVersion"

$about.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for V. I. Lenin (Kazakhstan) Coal Mine, Kazakhstan, M1438, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S ("build_version") holds the version string for each data row (rows 2-15).
$lastRow = $data.Cells.Item($data.Rows.Count, 19).End(-4162).Row  # -4162 = xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $data.Cells.Item($r, 19)  # column S
    if ($cell.Text -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
